$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.27203888085698
$ws.Range("C2").Value = 7.084572124863138
$ws.Range("D2").Value = 10.45597354189872
$ws.Range("F2").Value = 31.70197540192411
$ws.Range("G2").Value = 3.637569065551868
$ws.Range("J2").Value = 11.06776927751137
$ws.Range("M2").Value = 17.34413240756627
$ws.Range("N2").Value = 17.67072498597568
$ws.Range("O2").Value = 23.14078447566093
$ws.Range("B3").Value = 12.79903856225786
$ws.Range("C3").Value = 6.657861896246311
$ws.Range("D3").Value = 10.44229342651588
$ws.Range("F3").Value = 31.71141831923629
$ws.Range("G3").Value = 3.639749351830011
$ws.Range("J3").Value = 11.09396232882589
$ws.Range("M3").Value = 17.19908823864537
$ws.Range("N3").Value = 17.72789749605345
$ws.Range("O3").Value = 23.18266553470111
$ws.Range("B4").Value = 12.50147432924666
$ws.Range("C4").Value = 6.380264284868953
$ws.Range("D4").Value = 10.43544051920198
$ws.Range("F4").Value = 31.72587714722316
$ws.Range("G4").Value = 3.641159532637534
$ws.Range("J4").Value = 11.11158790185934
$ws.Range("M4").Value = 17.11236247483355
$ws.Range("N4").Value = 17.76484286231151
$ws.Range("O4").Value = 23.21437748494183
$ws.Range("B5").Value = 12.37860906608532
$ws.Range("C5").Value = 6.26324163770474
$ws.Range("D5").Value = 10.43303910207635
$ws.Range("F5").Value = 31.73394468762481
$ws.Range("G5").Value = 3.641752223775823
$ws.Range("J5").Value = 11.11915844335343
$ws.Range("M5").Value = 17.07763960728647
$ws.Range("N5").Value = 17.78036252317901
$ws.Range("O5").Value = 23.22880406061364
$ws.Range("B6").Value = 12.3581163252368
$ws.Range("C6").Value = 6.243575400836955
$ws.Range("D6").Value = 10.43266403760676
$ws.Range("F6").Value = 31.73541559200172
$ws.Range("G6").Value = 3.641851730414528
$ws.Range("J6").Value = 11.12043895852199
$ws.Range("M6").Value = 17.07191218216852
$ws.Range("N6").Value = 17.78296761664036
$ws.Range("O6").Value = 23.23129026786638
$ws.Range("B7").Value = 12.4998235613357
$ws.Range("C7").Value = 6.378701830174492
$ws.Range("D7").Value = 10.43540654608354
$ws.Range("F7").Value = 31.72597714478151
$ws.Range("G7").Value = 3.641167452793485
$ws.Range("J7").Value = 11.11168842994104
$ws.Range("M7").Value = 17.11189164374949
$ws.Range("N7").Value = 17.76505028515813
$ws.Range("O7").Value = 23.21456596440028
$ws.Range("B8").Value = 13.11053875302103
$ws.Range("C8").Value = 6.940690591669128
$ws.Range("D8").Value = 10.45093703071228
$ws.Range("F8").Value = 31.7034327270385
$ws.Range("G8").Value = 3.638306028177488
$ws.Range("J8").Value = 11.07648040549013
$ws.Range("M8").Value = 17.2936568185116
$ws.Range("N8").Value = 17.69005660650178
$ws.Range("O8").Value = 23.1539780242778
$ws.Range("B9").Value = 14.24405346480193
$ws.Range("C9").Value = 7.918353713140169
$ws.Range("D9").Value = 10.49355743269644
$ws.Range("F9").Value = 31.7279920504852
$ws.Range("G9").Value = 3.633259281782129
$ws.Range("J9").Value = 11.01968025325703
$ws.Range("M9").Value = 17.66713987210077
$ws.Range("N9").Value = 17.55755159666633
$ws.Range("O9").Value = 23.08291576116294
$ws.Range("B10").Value = 15.02902936457002
$ws.Range("C10").Value = 8.56030504557434
$ws.Range("D10").Value = 10.53213068749993
$ws.Range("F10").Value = 31.78794129653967
$ws.Range("G10").Value = 3.629891886679638
$ws.Range("J10").Value = 10.98541276212337
$ws.Range("M10").Value = 17.94991729532204
$ws.Range("N10").Value = 17.46900217109954
$ws.Range("O10").Value = 23.06001439535414
$ws.Range("B11").Value = 15.37423297662803
$ws.Range("C11").Value = 8.835740892264496
$ws.Range("D11").Value = 10.55121730727674
$ws.Range("F11").Value = 31.82427895689543
$ws.Range("G11").Value = 3.628433113190126
$ws.Range("J11").Value = 10.97144408153282
$ws.Range("M11").Value = 18.07994103085872
$ws.Range("N11").Value = 17.43061482775194
$ws.Range("O11").Value = 23.05599015872033
$ws.Range("B12").Value = 15.50313689205283
$ws.Range("C12").Value = 8.937656072533105
$ws.Range("D12").Value = 10.55866269813642
$ws.Range("F12").Value = 31.83933735710488
$ws.Range("G12").Value = 3.627891161729676
$ws.Range("J12").Value = 10.9663874025355
$ws.Range("M12").Value = 18.12933984078669
$ws.Range("N12").Value = 17.41634982827334
$ws.Range("O12").Value = 23.0553870717086
$ws.Range("B13").Value = 15.47545752437697
$ws.Range("C13").Value = 8.915812938108328
$ws.Range("D13").Value = 10.55704957833549
$ws.Range("F13").Value = 31.83603663296311
$ws.Range("G13").Value = 3.628007416553432
$ws.Range("J13").Value = 10.96746608759866
$ws.Range("M13").Value = 18.11869433538514
$ws.Range("N13").Value = 17.41940999277703
$ws.Range("O13").Value = 23.05547598912198
$ws.Range("B14").Value = 15.38487496933835
$ws.Range("C14").Value = 8.844173338471972
$ws.Range("D14").Value = 10.55182550276454
$ws.Range("F14").Value = 31.82549184201999
$ws.Range("G14").Value = 3.628388317299085
$ws.Range("J14").Value = 10.9710233965508
$ws.Range("M14").Value = 18.08400207745137
$ws.Range("N14").Value = 17.42943580258311
$ws.Range("O14").Value = 23.05592208247177
$ws.Range("B15").Value = 15.32915091168999
$ws.Range("C15").Value = 8.799981253765154
$ws.Range("D15").Value = 10.54865384433168
$ws.Range("F15").Value = 31.81920171309385
$ws.Range("G15").Value = 3.628622989903923
$ws.Range("J15").Value = 10.9732326896939
$ws.Range("M15").Value = 18.06277199733788
$ws.Range("N15").Value = 17.43561222459787
$ws.Range("O15").Value = 23.05631527349447
$ws.Range("B16").Value = 15.00622077436713
$ws.Range("C16").Value = 8.541970707429797
$ws.Range("D16").Value = 10.53091400099747
$ws.Range("F16").Value = 31.78574853902451
$ws.Range("G16").Value = 3.62998868678638
$ws.Range("J16").Value = 10.98635824564298
$ws.Range("M16").Value = 17.94144461578217
$ws.Range("N16").Value = 17.47154891073035
$ws.Range("O16").Value = 23.06040618096834
$ws.Range("B17").Value = 14.80498807500203
$ws.Range("C17").Value = 8.379439073151641
$ws.Range("D17").Value = 10.52042293048951
$ws.Range("F17").Value = 31.76754471784136
$ws.Range("G17").Value = 3.630845174963801
$ws.Range("J17").Value = 10.99482524222494
$ws.Range("M17").Value = 17.86734253093158
$ws.Range("N17").Value = 17.49407939074556
$ws.Range("O17").Value = 23.06455460611362
$ws.Range("B18").Value = 14.68813113684537
$ws.Range("C18").Value = 8.284393429126254
$ws.Range("D18").Value = 10.51453387823821
$ws.Range("F18").Value = 31.75792817244405
$ws.Range("G18").Value = 3.631344685292036
$ws.Range("J18").Value = 10.99984769423678
$ws.Range("M18").Value = 17.82485377785553
$ws.Range("N18").Value = 17.50721667565373
$ws.Range("O18").Value = 23.06754243473275
$ws.Range("B19").Value = 14.64837770811606
$ws.Range("C19").Value = 8.251944431171927
$ws.Range("D19").Value = 10.51256498150274
$ws.Range("F19").Value = 31.75481896469006
$ws.Range("G19").Value = 3.63151499434889
$ws.Range("J19").Value = 11.0015743939182
$ws.Range("M19").Value = 17.81049177677127
$ws.Range("N19").Value = 17.51169539500732
$ws.Range("O19").Value = 23.06865736239604
$ws.Range("B20").Value = 14.82652567856178
$ws.Range("C20").Value = 8.396902489371492
$ws.Range("D20").Value = 10.52152472792753
$ws.Range("F20").Value = 31.76939421530695
$ws.Range("G20").Value = 3.630753288565622
$ws.Range("J20").Value = 10.99390813668455
$ws.Range("M20").Value = 17.87521735313195
$ws.Range("N20").Value = 17.49166253108827
$ws.Range("O20").Value = 23.064050706588
$ws.Range("B21").Value = 15.41153137354508
$ws.Range("C21").Value = 8.865280386031605
$ws.Range("D21").Value = 10.55335406232126
$ws.Range("F21").Value = 31.82855392464025
$ws.Range("G21").Value = 3.628276154209082
$ws.Range("J21").Value = 10.96997220608337
$ws.Range("M21").Value = 18.09418794861271
$ws.Range("N21").Value = 17.42648361935548
$ws.Range("O21").Value = 23.05576605585183
$ws.Range("B22").Value = 15.78323358720843
$ws.Range("C22").Value = 9.157483794040401
$ws.Range("D22").Value = 10.57542345581773
$ws.Range("F22").Value = 31.87478094495084
$ws.Range("G22").Value = 3.626718115562094
$ws.Range("J22").Value = 10.95568654978669
$ws.Range("M22").Value = 18.2382226009439
$ws.Range("N22").Value = 17.38546725860113
$ws.Range("O22").Value = 23.05571888385917
$ws.Range("B23").Value = 15.58585426490859
$ws.Range("C23").Value = 9.002801651626932
$ws.Range("D23").Value = 10.56352993403096
$ws.Range("F23").Value = 31.84941897780399
$ws.Range("G23").Value = 3.627544114137567
$ws.Range("J23").Value = 10.96318682159911
$ws.Range("M23").Value = 18.16127645304474
$ws.Range("N23").Value = 17.40721402799085
$ws.Range("O23").Value = 23.05525266095609
$ws.Range("B24").Value = 14.81679215014772
$ws.Range("C24").Value = 8.389012275120708
$ws.Range("D24").Value = 10.52102616146408
$ws.Range("F24").Value = 31.76855541224856
$ws.Range("G24").Value = 3.630794808273735
$ws.Range("J24").Value = 10.99432227833613
$ws.Range("M24").Value = 17.87165678951249
$ws.Range("N24").Value = 17.49275461943854
$ws.Range("O24").Value = 23.06427664181916
$ws.Range("B25").Value = 13.94523756941284
$ws.Range("C25").Value = 7.667255372727612
$ws.Range("D25").Value = 10.48074033843802
$ws.Range("F25").Value = 31.71398193021379
$ws.Range("G25").Value = 3.634564509281321
$ws.Range("J25").Value = 11.03373547464611
$ws.Range("M25").Value = 17.56448338534523
$ws.Range("N25").Value = 17.59184682850477
$ws.Range("O25").Value = 23.09700517200572
